$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell holding the plain default style (no explicit formatting),
# used to restore style on cells where we must force text via a leading
# apostrophe (to stop Excel from re-interpreting numeric-looking text as a number).
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "79.532.03"
$ws.Range("E2").Value = "  +3.96%  "
$ws.Range("D3").Value = "3.206.10"
$ws.Range("E3").Value = "  +7.10%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'207.39"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").Value = "'630.70"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.228"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +12.45%  "
$ws.Range("D9").Value = "'0.583"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  +6.23%  "
$ws.Range("D10").Value = "3.205.39"
$ws.Range("E10").Value = "  +7.22%  "
$ws.Range("D11").Value = "'0.584"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +35.10%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  +8.22%  "
$ws.Range("D14").Value = "3.794.68"
$ws.Range("E14").Value = "  +7.28%  "
$ws.Range("D15").Value = "'0.0000228"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +21.27%  "
$ws.Range("D16").Value = "'31.83"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +9.10%  "
$ws.Range("D17").Value = "79.099.20"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "3.201.59"
$ws.Range("E18").Value = "  +7.31%  "
$ws.Range("D19").Value = "'14.49"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +7.61%  "
$ws.Range("D20").Value = "'9.47"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +5.94%  "
$ws.Range("D21").Value = "'433.48"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +16.12%  "
$ws.Range("D22").Value = "'2.85"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +25.62%  "
$ws.Range("D23").Value = "'4.99"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +16.14%  "
$ws.Range("D24").Value = "'6.82"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +6.14%  "
$ws.Range("D25").Value = "3.371.27"
$ws.Range("E25").Value = "  +7.58%  "
$ws.Range("D26").Value = "'4.78"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +9.78%  "
$ws.Range("D27").Value = "'77.17"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +6.01%  "
$ws.Range("D28").Value = "'11.03"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +12.55%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +8.19%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +7.82%  "
$ws.Range("D33").Value = "'1.47"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +6.05%  "
$ws.Range("D34").Value = "'518.92"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("D36").Value = "'0.137"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +23.64%  "
$ws.Range("D37").Value = "'22.82"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +12.00%  "
$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D38").Value = "'0.121"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +14.94%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.406"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +6.53%  "
$ws.Range("D41").Value = "'164.44"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "'197.21"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +5.51%  "
$ws.Range("D43").Value = "'20.01"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "'5.45"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +10.00%  "
$ws.Range("E46").Value = "  +14.23%  "
$ws.Range("E47").Value = "  +8.62%  "
$ws.Range("E48").Value = "  +5.69%  "
$ws.Range("D49").Value = "'43.14"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").Value = "'2.54"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +9.69%  "
$ws.Range("D51").Value = "'0.629"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +6.66%  "
